{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change (from the diff): the final paragraph that currently reads\n// \"\\tAlternatively\" gets grammar-check proof marks wrapped around\n// \"Alternatively\" and a large chunk of additional sentences appended after\n// it, all still inside the same paragraph (the trailing _GoBack bookmark\n// stays at the very end of the paragraph).\n\n// Locate the (single) \"Alternatively\" run in the document body.\nconst results = context.document.body.search(\"Alternatively\", { matchCase: true, matchWholeWord: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Alternatively\" in the document body.');\n}\n\nconst target = results.items[0];\n\n// Flat-OPC wrapped WordprocessingML fragment that replaces the \"Alternatively\"\n// range with: gramStart proofErr, the \"Alternatively\" run, gramEnd proofErr,\n// and then the newly authored sentences as additional runs.\nconst replacementOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r><w:t>Alternatively</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">after registration is complete, the student can view the current active courses in their schedule.  This </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">is relatively common for students to use CURSE for </w:t></w:r>' +\n  '<w:r><w:t>day to day schedule viewing.</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">  Viewing the professors provides a place to contact </w:t></w:r>' +\n  '<w:r><w:t>instructors and request assistance.</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word.ActiveDocument is the open document.\n#\n# Target change (from the diff): the final paragraph that currently reads\n# \"\\tAlternatively\" gets grammar-check proof marks wrapped around\n# \"Alternatively\" and a large chunk of additional sentences appended after\n# it, all still inside the same paragraph (the trailing _GoBack bookmark\n# stays at the very end of the paragraph).\n\n$d = $word.ActiveDocument\n\n# Locate the (single) \"Alternatively\" run in the document.\n$found = $d.Content\n$found.Find.ClearFormatting()\n$found.Find.Execute(\"Alternatively\")\n\n# Re-seat into a plain Range over the same span so the edit below lands\n# exactly on \"Alternatively\" (not the live Find range).\n$target = $d.Range($found.Start, $found.End)\n\n# Flat-OPC wrapped WordprocessingML fragment that replaces the \"Alternatively\"\n# range with: gramStart proofErr, the \"Alternatively\" run, gramEnd proofErr,\n# and then the newly authored sentences as additional runs.\n$replacementOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>Alternatively</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">after registration is complete, the student can view the current active courses in their schedule.  This </w:t></w:r><w:r><w:t xml:space=\"preserve\">is relatively common for students to use CURSE for </w:t></w:r><w:r><w:t>day to day schedule viewing.</w:t></w:r><w:r><w:t xml:space=\"preserve\">  Viewing the professors provides a place to contact </w:t></w:r><w:r><w:t>instructors and request assistance.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$target.InsertXML($replacementOoxml)\n"}
